$wb = $excel.ActiveWorkbook

# ALC!row 106
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 1923.4231
$ws.Range("I106").Value = 1300.6923
$ws.Range("K106").Value = 1300.6923
$ws.Range("M106").Value = -669.6922999999999

# ARM!row 33
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 15280.556
$ws.Range("I33").Value = 10841.667
$ws.Range("J33").Value = 17500
$ws.Range("K33").Value = 10841.667
$ws.Range("L33").Value = 17500
$ws.Range("M33").Value = -10512.667
$ws.Range("N33").Value = -18158

# ARM!row 80
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 26125.111
$ws.Range("J80").Value = 26125.111
$ws.Range("L80").Value = 26125.111
$ws.Range("N80").Value = -28121.111

# ARM!row 83
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H83").Value = 26125.111
$ws.Range("J83").Value = 26125.111
$ws.Range("L83").Value = 78375.333
$ws.Range("N83").Value = -88359.333

# ARM!row 110
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 1663.1852
$ws.Range("I110").Value = 618.1053000000001
$ws.Range("J110").Value = 4145.25
$ws.Range("K110").Value = 618.1053000000001
$ws.Range("L110").Value = 4145.25
$ws.Range("M110").Value = 1426.8947
$ws.Range("N110").Value = -8235.25

# ARM!row 113
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H113").Value = 29750
$ws.Range("J113").Value = 29750
$ws.Range("L113").Value = 29750
$ws.Range("N113").Value = -38428

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1827.6666
$ws.Range("I132").Value = 1320.5098
$ws.Range("J132").Value = 3983.0833
$ws.Range("K132").Value = 3961.5294
$ws.Range("L132").Value = 11949.2499
$ws.Range("M132").Value = -1431.5294
$ws.Range("N132").Value = -17009.2499

# BSM!row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1727.9286
$ws.Range("I105").Value = 1640
$ws.Range("J105").Value = 2255.5
$ws.Range("K105").Value = 1640
$ws.Range("L105").Value = 2255.5
$ws.Range("M105").Value = 107
$ws.Range("N105").Value = -5749.5

# CRP!row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 21905.75
$ws.Range("I10").Value = 606.125
$ws.Range("K10").Value = 606.125
$ws.Range("M10").Value = -467.125

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2040.8572
$ws.Range("I31").Value = 1275.711
$ws.Range("J31").Value = 2690.5095
$ws.Range("K31").Value = 1275.711
$ws.Range("L31").Value = 2690.5095
$ws.Range("M31").Value = -980.711
$ws.Range("N31").Value = -3280.5095

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2040.8572
$ws.Range("I34").Value = 1275.711
$ws.Range("J34").Value = 2690.5095
$ws.Range("K34").Value = 1275.711
$ws.Range("L34").Value = 2690.5095
$ws.Range("M34").Value = -1073.711
$ws.Range("N34").Value = -3094.5095

# CUL!row 2
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 28.38889
$ws.Range("I2").Value = 17
$ws.Range("J2").Value = 42.625
$ws.Range("K2").Value = 102
$ws.Range("L2").Value = 255.75
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = -481.75

# CUL!row 18
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 703.4545000000001
$ws.Range("I18").Value = 248.42857
$ws.Range("J18").Value = 1499.75
$ws.Range("K18").Value = 745.28571
$ws.Range("L18").Value = 4499.25
$ws.Range("M18").Value = -576.28571
$ws.Range("N18").Value = -4837.25

# CUL!row 75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M75").ClearContents()
$ws.Range("H75").Value = 2999
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2999
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 8997
$ws.Range("N75").Value = -10993

# CUL!row 78
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M78").ClearContents()
$ws.Range("H78").Value = 2999
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2999
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 26991
$ws.Range("N78").Value = -36975

# CUL!row 116
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1382.375
$ws.Range("I116").Value = 171.8
$ws.Range("K116").Value = 515.4000000000001
$ws.Range("M116").Value = 2926.6

# CUL!row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2278.2
$ws.Range("I136").Value = 1578.8889
$ws.Range("J136").Value = 3327.1667
$ws.Range("K136").Value = 4736.6667
$ws.Range("L136").Value = 9981.500100000001
$ws.Range("M136").Value = 363.3333000000002
$ws.Range("N136").Value = -20181.5001

# CUL!row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 3223.6316
$ws.Range("I137").Value = 2273.25
$ws.Range("J137").Value = 4852.857
$ws.Range("K137").Value = 6819.75
$ws.Range("L137").Value = 14558.571
$ws.Range("M137").Value = -1719.75
$ws.Range("N137").Value = -24758.571

# GSM!row 11
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 3484999.8
$ws.Range("I11").Value = 4314687
$ws.Range("J11").Value = 2463846.2
$ws.Range("K11").Value = 4314687
$ws.Range("L11").Value = 2463846.2
$ws.Range("M11").Value = -4314548
$ws.Range("N11").Value = -2464124.2

# GSM!row 15
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 20965.834
$ws.Range("I15").Value = 9996
$ws.Range("J15").Value = 23159.8
$ws.Range("K15").Value = 9996
$ws.Range("L15").Value = 23159.8
$ws.Range("M15").Value = -9708
$ws.Range("N15").Value = -23735.8

# GSM!row 20
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1366668.6
$ws.Range("I20").Value = 2000000
$ws.Range("K20").Value = 2000000
$ws.Range("M20").Value = -1999755

# GSM!row 81
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H81").Value = 20965.834
$ws.Range("I81").Value = 9996
$ws.Range("J81").Value = 23159.8
$ws.Range("K81").Value = 9996
$ws.Range("L81").Value = 23159.8
$ws.Range("M81").Value = -8998
$ws.Range("N81").Value = -25155.8

# GSM!row 84
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H84").Value = 20965.834
$ws.Range("I84").Value = 9996
$ws.Range("J84").Value = 23159.8
$ws.Range("K84").Value = 29988
$ws.Range("L84").Value = 69479.39999999999
$ws.Range("M84").Value = -24996
$ws.Range("N84").Value = -79463.39999999999

# GSM!row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1532.409
$ws.Range("I97").Value = 1142.7894
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 1142.7894
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = -646.7893999999999
$ws.Range("N97").Value = -4992

# GSM!row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 46941.523
$ws.Range("I102").Value = 2609.5
$ws.Range("K102").Value = 2609.5
$ws.Range("M102").Value = -987.5

# LTW!row 24
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 18434.143
$ws.Range("J24").Value = 18434.143
$ws.Range("L24").Value = 18434.143
$ws.Range("N24").Value = -19120.143

# LTW!row 35
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 10134
$ws.Range("J35").Value = 20000
$ws.Range("L35").Value = 20000
$ws.Range("N35").Value = -20672

# LTW!row 93
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2096.8823
$ws.Range("I93").Value = 1624.2727
$ws.Range("J93").Value = 2963.3333
$ws.Range("K93").Value = 1624.2727
$ws.Range("L93").Value = 2963.3333
$ws.Range("M93").Value = -376.2727
$ws.Range("N93").Value = -5459.3333

# LTW!row 106
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 28857.143
$ws.Range("J106").Value = 28857.143
$ws.Range("L106").Value = 28857.143
$ws.Range("N106").Value = -31381.143

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2109.0356
$ws.Range("I132").Value = 1406.7368
$ws.Range("J132").Value = 3591.6667
$ws.Range("K132").Value = 4220.2104
$ws.Range("L132").Value = 10775.0001
$ws.Range("M132").Value = -1690.2104
$ws.Range("N132").Value = -15835.0001

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4436.2666
$ws.Range("I136").Value = 3662.2632
$ws.Range("J136").Value = 5773.1816
$ws.Range("K136").Value = 10986.7896
$ws.Range("L136").Value = 17319.5448
$ws.Range("M136").Value = -8436.7896
$ws.Range("N136").Value = -22419.5448

# WVR!row 13
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M13").ClearContents()
$ws.Range("H13").Value = 67670.664
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 67670.664
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 67670.664
$ws.Range("N13").Value = -67950.664

# WVR!row 98
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H98").Value = 29949.75
$ws.Range("J98").Value = 29949.75
$ws.Range("L98").Value = 29949.75
$ws.Range("N98").Value = -35939.75

# WVR!row 104
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 37000
$ws.Range("J104").Value = 37000
$ws.Range("L104").Value = 37000
$ws.Range("N104").Value = -43988

# WVR!row 119
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H119").Value = 29137.143
$ws.Range("J119").Value = 29137.143
$ws.Range("L119").Value = 29137.143
$ws.Range("N119").Value = -38813.143
